$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.629231666666667
$ws.Range("H2").Value = 7.887695
$ws.Range("I2").Value = 0.1414315557047068
$ws.Range("J2").Value = 0.1414315557047067
$ws.Range("M2").Value = 0.9956583333333334
$ws.Range("N2").Value = 2.986975
$ws.Range("O2").Value = 0.1535710477437721
$ws.Range("P2").Value = 0.153571047743772
$ws.Range("Q2").Value = 2.617816419180556
$ws.Range("R2").Value = 23.560347772625
$ws.Range("S2").Value = 0.02171979219360348
$ws.Range("T2").Value = 0.02171979219360347

# Row 3
$ws.Range("G3").Value = 2.629231666666667
$ws.Range("H3").Value = 7.887695
$ws.Range("I3").Value = 0.1414315557047068
$ws.Range("J3").Value = 0.1414315557047067
$ws.Range("O3").Value = 0.2875834480798523
$ws.Range("P3").Value = 0.2875834480798522
$ws.Range("Q3").Value = 4.902230487637778
$ws.Range("R3").Value = 44.12007438873999
$ws.Range("S3").Value = 0.04067337445685727
$ws.Range("T3").Value = 0.04067337445685725

# Row 4
$ws.Range("G4").Value = 2.629231666666667
$ws.Range("H4").Value = 7.887695
$ws.Range("I4").Value = 0.1414315557047068
$ws.Range("J4").Value = 0.1414315557047067
$ws.Range("M4").Value = 2.910118
$ws.Range("N4").Value = 8.730354
$ws.Range("O4").Value = 0.4488586650219809
$ws.Range("P4").Value = 0.4488586650219808
$ws.Range("Q4").Value = 7.651374399336667
$ws.Range("R4").Value = 68.86236959403
$ws.Range("S4").Value = 0.0634827792855966
$ws.Range("T4").Value = 0.06348277928559658

# Row 5
$ws.Range("G5").Value = 2.629231666666667
$ws.Range("H5").Value = 7.887695
$ws.Range("I5").Value = 0.1414315557047068
$ws.Range("J5").Value = 0.1414315557047067
$ws.Range("M5").Value = 0.7130856666666667
$ws.Range("N5").Value = 2.139257
$ws.Range("O5").Value = 0.1099868391543949
$ws.Range("P5").Value = 0.1099868391543948
$ws.Range("Q5").Value = 1.874867415846111
$ws.Range("R5").Value = 16.873806742615
$ws.Range("S5").Value = 0.01555560976864942
$ws.Range("T5").Value = 0.01555560976864941

# Row 6
$ws.Range("I6").Value = 0.6147160060020365
$ws.Range("J6").Value = 0.6147160060020365
$ws.Range("M6").Value = 0.9956583333333334
$ws.Range("N6").Value = 2.986975
$ws.Range("O6").Value = 0.1535710477437721
$ws.Range("P6").Value = 0.153571047743772
$ws.Range("Q6").Value = 11.37803827177778
$ws.Range("R6").Value = 102.402344446
$ws.Range("S6").Value = 0.09440258110659962
$ws.Range("T6").Value = 0.0944025811065996

# Row 7
$ws.Range("I7").Value = 0.6147160060020365
$ws.Range("J7").Value = 0.6147160060020365
$ws.Range("O7").Value = 0.2875834480798523
$ws.Range("P7").Value = 0.2875834480798522
$ws.Range("S7").Value = 0.1767821485959408
$ws.Range("T7").Value = 0.1767821485959408

# Row 8
$ws.Range("I8").Value = 0.6147160060020365
$ws.Range("J8").Value = 0.6147160060020365
$ws.Range("M8").Value = 2.910118
$ws.Range("N8").Value = 8.730354
$ws.Range("O8").Value = 0.4488586650219809
$ws.Range("P8").Value = 0.4488586650219808
$ws.Range("Q8").Value = 33.25581966309333
$ws.Range("R8").Value = 299.30237696784
$ws.Range("S8").Value = 0.2759206058217181
$ws.Range("T8").Value = 0.275920605821718

# Row 9
$ws.Range("I9").Value = 0.6147160060020365
$ws.Range("J9").Value = 0.6147160060020365
$ws.Range("M9").Value = 0.7130856666666667
$ws.Range("N9").Value = 2.139257
$ws.Range("O9").Value = 0.1099868391543949
$ws.Range("P9").Value = 0.1099868391543948
$ws.Range("Q9").Value = 8.148895795635557
$ws.Range("R9").Value = 73.34006216072001
$ws.Range("S9").Value = 0.06761067047777802
$ws.Range("T9").Value = 0.067610670477778

# Row 10
$ws.Range("G10").Value = 4.24731
$ws.Range("H10").Value = 12.74193
$ws.Range("I10").Value = 0.2284711798035388
$ws.Range("J10").Value = 0.2284711798035388
$ws.Range("M10").Value = 0.9956583333333334
$ws.Range("N10").Value = 2.986975
$ws.Range("O10").Value = 0.1535710477437721
$ws.Range("P10").Value = 0.153571047743772
$ws.Range("Q10").Value = 4.22886959575
$ws.Range("R10").Value = 38.05982636175001
$ws.Range("S10").Value = 0.03508655846168519
$ws.Range("T10").Value = 0.03508655846168518

# Row 11
$ws.Range("G11").Value = 4.24731
$ws.Range("H11").Value = 12.74193
$ws.Range("I11").Value = 0.2284711798035388
$ws.Range("J11").Value = 0.2284711798035388
$ws.Range("O11").Value = 0.2875834480798523
$ws.Range("P11").Value = 0.2875834480798522
$ws.Range("Q11").Value = 7.919154799639999
$ws.Range("R11").Value = 71.27239319675999
$ws.Range("S11").Value = 0.0657045296747736
$ws.Range("T11").Value = 0.06570452967477357

# Row 12
$ws.Range("G12").Value = 4.24731
$ws.Range("H12").Value = 12.74193
$ws.Range("I12").Value = 0.2284711798035388
$ws.Range("J12").Value = 0.2284711798035388
$ws.Range("M12").Value = 2.910118
$ws.Range("N12").Value = 8.730354
$ws.Range("O12").Value = 0.4488586650219809
$ws.Range("P12").Value = 0.4488586650219808
$ws.Range("Q12").Value = 12.36017328258
$ws.Range("R12").Value = 111.24155954322
$ws.Range("S12").Value = 0.1025512687626134
$ws.Range("T12").Value = 0.1025512687626134

# Row 13
$ws.Range("G13").Value = 4.24731
$ws.Range("H13").Value = 12.74193
$ws.Range("I13").Value = 0.2284711798035388
$ws.Range("J13").Value = 0.2284711798035388
$ws.Range("M13").Value = 0.7130856666666667
$ws.Range("N13").Value = 2.139257
$ws.Range("O13").Value = 0.1099868391543949
$ws.Range("P13").Value = 0.1099868391543948
$ws.Range("Q13").Value = 3.02869588289
$ws.Range("R13").Value = 27.25826294601
$ws.Range("S13").Value = 0.02512882290446665
$ws.Range("T13").Value = 0.02512882290446664

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2859396666666667
$ws.Range("H14").Value = 0.8578190000000001
$ws.Range("I14").Value = 0.01538125848971795
$ws.Range("J14").Value = 0.01538125848971795
$ws.Range("M14").Value = 0.9956583333333334
$ws.Range("N14").Value = 2.986975
$ws.Range("O14").Value = 0.1535710477437721
$ws.Range("P14").Value = 0.153571047743772
$ws.Range("Q14").Value = 0.2846982119472223
$ws.Range("R14").Value = 2.562283907525
$ws.Range("S14").Value = 0.002362115981883775
$ws.Range("T14").Value = 0.002362115981883775

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2859396666666667
$ws.Range("H15").Value = 0.8578190000000001
$ws.Range("I15").Value = 0.01538125848971795
$ws.Range("J15").Value = 0.01538125848971795
$ws.Range("O15").Value = 0.2875834480798523
$ws.Range("P15").Value = 0.2875834480798522
$ws.Range("Q15").Value = 0.5331375585231112
$ws.Range("R15").Value = 4.798238026708001
$ws.Range("S15").Value = 0.00442339535228059
$ws.Range("T15").Value = 0.004423395352280588

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2859396666666667
$ws.Range("H16").Value = 0.8578190000000001
$ws.Range("I16").Value = 0.01538125848971795
$ws.Range("J16").Value = 0.01538125848971795
$ws.Range("M16").Value = 2.910118
$ws.Range("N16").Value = 8.730354
$ws.Range("O16").Value = 0.4488586650219809
$ws.Range("P16").Value = 0.4488586650219808
$ws.Range("Q16").Value = 0.8321181708806669
$ws.Range("R16").Value = 7.489063537926001
$ws.Range("S16").Value = 0.00690401115205281
$ws.Range("T16").Value = 0.006904011152052809

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2859396666666667
$ws.Range("H17").Value = 0.8578190000000001
$ws.Range("I17").Value = 0.01538125848971795
$ws.Range("J17").Value = 0.01538125848971795
$ws.Range("M17").Value = 0.7130856666666667
$ws.Range("N17").Value = 2.139257
$ws.Range("O17").Value = 0.1099868391543949
$ws.Range("P17").Value = 0.1099868391543948
$ws.Range("Q17").Value = 0.2038994778314445
$ws.Range("R17").Value = 1.835095300483
$ws.Range("S17").Value = 0.001691736003500779
$ws.Range("T17").Value = 0.001691736003500778
